$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (columns A..AH), replacing the old readings with the updated dataset.
$newData = @(
    @(45045.50694444445, 10.232, 6.988, 3.382, 22.493, 16.551, 7.711, 22.993, 12.747, 5.033, 6.872, 8.682, 9.968, 2.214, 8.271000000000001, 11.091, 7.662, 2.404, 1.071, 118.497, 22.775, 7.635, 14.429, 7.663, 2.148, 12.965, 6.744, 6.372, 7.242, 9.571, 2.474, 20.428, 3.858, 9.544),
    @(45045.51388888889, 22.356, 16.537, 1.957, 49.025, 39.408, 17.406, 64.995, 27.294, 12.061, 17.391, 19.562, 21.073, 5.412, 17.67, 24.915, 15.179, 1.301, 0.9409999999999999, 261.641, 49.364, 16.31, 32.893, 17.279, 2.727, 32.929, 14.407, 12.915, 15.127, 20.688, 1.101, 59.249, 9.018000000000001, 20.392),
    @(45045.52083333334, 14.245, 10.559, 1.231, 31.314, 25.129, 11.077, 46.248, 17.408, 7.751, 11.042, 12.512, 13.482, 3.436, 11.279, 15.901, 9.731999999999999, 0.873, 0.576, 164.364, 31.616, 10.411, 21.03, 11.022, 1.732, 22.511, 9.196, 8.263, 9.682, 13.211, 0.708, 42.186, 5.739, 13.016),
    @(45045.52777777778, 21.01, 15.69, 1.25, 45.96, 37.43, 16.42, 63.8, 25.55, 11.46, 16.66, 18.41, 19.66, 5.17, 16.54, 23.47, 14, 0.75, 0.72, 244.48, 46.24, 15.27, 31.03, 16.29, 2.35, 31.37, 13.49, 11.98, 14.07, 19.39, 0.52, 57.89, 8.56, 19.09)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}

# Row 6 (old extra data row) is removed entirely; sheet now spans only through row 5.
$ws.Rows.Item(6).Delete()

# Widen several numeric columns by one unit (7->8, and column T/20 8->9) to fit the new, larger values.
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666

